# "add flow statisitics in wdcp"
#
# Inserts a new row into the request_type table on Sheet1 documenting the
# new REQ_TYPE_FLOW_STATISTICS request (value 0x03), right after the
# existing REQ_TYPE_FAKE_AP row. Inserting the row pushes the
# encrypt_type table (previously rows 19-22) down to rows 20-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a blank row at 18 - everything from row 18 downward (the
# encrypt_type section) shifts down by one row.
$ws.Rows.Item(18).Insert()

# Fill in the new request_type/value/meaning entry in the freshly
# inserted row 18.
$ws.Range("A18").Value = "REQ_TYPE_FLOW_STATISTICS"
$ws.Range("B18").Value = "0x03"
$ws.Range("C18").Value = "获取当前设备统计的数据流量信息"

# Match the author's final selection state.
$ws.Range("C23").Select()
